$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")
$ws.Activate()

# New polling update: "Latest Morgan" (row 6) is pushed down to "Second Morgan"
# (row 7), the old "Second Morgan" is pushed down to "Third Morgan" (row 8,
# discarding its previous numbers), and the fresh Morgan poll figures are
# entered into row 6.
$ws.Range("B8:G8").Value2 = $ws.Range("B7:G7").Value2
$ws.Range("B7:G7").Value2 = $ws.Range("B6:G6").Value2

$ws.Range("B6").Value2 = 55.5
$ws.Range("C6").Value2 = 56
$ws.Range("D6").Value2 = 63.5
$ws.Range("E6").Value2 = 43.5
$ws.Range("F6").Value2 = 49
$ws.Range("G6").Value2 = 62.5

$ws.Range("I26:I27").Select()
